$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 10518
$wsExhibit.Range("G2").Value = 80
$wsExhibit.Range("F8").Value = 482
$wsExhibit.Range("F9").Value = 784
$wsExhibit.Range("F12").Value = 1112
$wsExhibit.Range("F13").Value = 3263
$wsExhibit.Range("F14").Value = 2425
$wsExhibit.Range("F16").Value = 2183
$wsExhibit.Range("F17").Value = 2183
$wsExhibit.Range("F22").Value = 590
$wsExhibit.Range("F23").Value = 70
$wsExhibit.Range("F25").Value = 12
$wsExhibit.Range("F26").Value = 30
$wsExhibit.Range("F32").Value = 410
$wsExhibit.Range("F33").Value = 611
$wsExhibit.Range("F34").Value = 33
$wsExhibit.Range("F35").Value = 58
$wsExhibit.Range("F36").Value = 275
$wsExhibit.Range("F37").Value = 13
$wsExhibit.Range("F38").Value = 1581
$wsExhibit.Range("F39").Value = 518
$wsExhibit.Range("F40").Value = 488
$wsExhibit.Range("F41").Value = 1736
$wsExhibit.Range("F42").Value = 148
$wsExhibit.Range("F43").Value = 452
$wsExhibit.Range("F45").Value = 477
$wsExhibit.Range("F46").Value = 1057
$wsExhibit.Range("F48").Value = 370

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 49

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 10518
$wsAll.Range("G2").Value = 80
$wsAll.Range("F10").Value = 482
$wsAll.Range("F11").Value = 784
$wsAll.Range("F12").Value = 1112
$wsAll.Range("F13").Value = 3263
$wsAll.Range("F14").Value = 2425
$wsAll.Range("F15").Value = 2183
$wsAll.Range("F16").Value = 2183
$wsAll.Range("F18").Value = 590
$wsAll.Range("F19").Value = 70
$wsAll.Range("F21").Value = 12
$wsAll.Range("F22").Value = 30
$wsAll.Range("F28").Value = 410
$wsAll.Range("F29").Value = 611
$wsAll.Range("F30").Value = 33
$wsAll.Range("F31").Value = 49
$wsAll.Range("F34").Value = 58
$wsAll.Range("F35").Value = 275
$wsAll.Range("F36").Value = 1581
$wsAll.Range("F37").Value = 518
$wsAll.Range("F39").Value = 488
$wsAll.Range("F40").Value = 1736
$wsAll.Range("F41").Value = 148
$wsAll.Range("F45").Value = 452
$wsAll.Range("F47").Value = 477
$wsAll.Range("F48").Value = 1057
$wsAll.Range("F49").Value = 370
